$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column range to text format so that numeric-looking
# strings (e.g. "317.02") are not auto-converted to numbers by Excel,
# matching the original inline-string cell type. Revert formatting
# afterwards so the cell style index is left unchanged.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.632.24'
$ws.Range('D3').Value = '2.517.12'
$ws.Range('D5').Value = '317.02'
$ws.Range('D6').Value = '95.51'
$ws.Range('D7').Value = '0.578'
$ws.Range('D9').Value = '0.536'
$ws.Range('D10').Value = '36.05'
$ws.Range('D11').Value = '0.0809'
$ws.Range('D12').Value = '7.69'
$ws.Range('D14').Value = '2.899.96'
$ws.Range('D15').Value = '15.46'
$ws.Range('D16').Value = '2.496.06'
$ws.Range('D17').Value = '0.860'
$ws.Range('D18').Value = '42.666.92'
$ws.Range('D19').Value = '13.11'
$ws.Range('D20').Value = '0.0₃0968'
$ws.Range('D22').Value = '71.22'
$ws.Range('D23').Value = '251.03'
$ws.Range('D25').Value = '2.03'
$ws.Range('D26').Value = '26.74'
$ws.Range('D29').Value = '38.70'
$ws.Range('D30').Value = '10.04'
$ws.Range('D31').Value = '5.90'
$ws.Range('D32').Value = '155.99'
$ws.Range('D33').Value = '19.73'
$ws.Range('D34').Value = '3.34'
$ws.Range('D35').Value = '2.08'
$ws.Range('D40').Value = '23.96'
$ws.Range('D45').Value = '2.056.93'
$ws.Range('D46').Value = '0.0300'
$ws.Range('D47').Value = '84.40'
$ws.Range('D48').Value = '8.80'
$ws.Range('D49').Value = '2.754.49'
$ws.Range('D50').Value = '73.07'
$ws.Range('D51').Value = '0.190'

$ws.Range("D2:D51").ClearFormats()

# Coin name / link (B, C) and Volume(1h) (E) columns are plain text and
# are not subject to numeric auto-conversion.
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E6').Value = '  -2.93%  '
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('E15').Value = '  +4.05%  '
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('E19').Value = '  -4.18%  '
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('E25').Value = '  -2.90%  '
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  +12.43%  '
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('E32').Value = '  -1.91%  '
$ws.Range('E33').Value = '  +4.41%  '
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  -3.97%  '
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('E37').Value = '  -5.21%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  -7.72%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E51').Value = '  -0.70%  '
